$d = $word.ActiveDocument

# Locate the paragraph that currently contains only "Funciones:"
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Funciones:`r") {
        $target = $p
        break
    }
}

# Paragraph formatting: spacing before/after 12pt (240 twentieths-of-a-point)
# with auto line rule, and right indent of 30pt (600 twentieths-of-a-point).
$target.Format.SpaceBefore = 12
$target.Format.SpaceAfter = 12
$target.Format.LineSpacingRule = 5
$target.Format.LeftIndent = 0
$target.Format.RightIndent = 30
$target.Format.FirstLineIndent = 0

# Make the existing "Funciones:" run bold.
$boldRange = $d.Range($target.Range.Start, $target.Range.End - 1)
$boldRange.Font.Bold = 1

# Append the new (non-bold) explanatory sentence as its own run.
$insertPos = $target.Range.End - 1
$appendRange = $d.Range($insertPos, $insertPos)
$appendRange.InsertAfter(" Se ha creado una carpeta de procesos que alberga programas auxiliares para el desarrollo del código principal. En particular, contiene dos programas dedicados a desencriptar ejemplos específicos del dataset de desarrollo. El objetivo de estos es identificar el formato de compresión utilizado, centrándose en los algoritmos RLE y LZ78.")
$newTextRange = $d.Range($insertPos, $target.Range.End - 1)
$newTextRange.Font.Bold = 0

# Insert a new empty paragraph right after, duplicating the bold paragraph
# mark / empty run formatting of the paragraph that used to follow "Funciones:".
$following = $d.Paragraphs.Item($target.Index + 1)
$following.Range.InsertParagraphBefore()
